# feat(core-lib): 新增 Guard 方法。
# Adds a new row (row 6) to the ResX resource table describing the new
# "ArgumentNullException_with_method_and_parameter_name" guard resource,
# mirroring the structure/format of the existing rows (2-5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (cell styles) of the last existing data row (row 5)
# down into the new row 6, so the new row's columns line up with the same
# look (A/B locked "label" style, C/E/F/G/H regular style) as every other
# data row.
$ws.Range("A5:H5").Copy()
$ws.Range("A6:H6").PasteSpecial(-4122)

# Column D ("Comment") is intentionally left blank for every data row, so
# drop the pasted placeholder cell in D6 to match.
$ws.Range("D6").Clear()

# Project / File columns repeat the same values as every other row.
$ws.Range("A6").Value = "CoreLib, Private.CoreLib"
$ws.Range("B6").Value = "Strings"

# New resource entry: key + the four localized message variants.
$ws.Range("C6").Value = "ArgumentNullException_with_method_and_parameter_name"
$ws.Range("E6").Value = "The method `"{0}`" may be missing the required parameter value: {1}."
$ws.Range("F6").Value = "The method `"{0}`" may be missing the required parameter value: {1}."
$ws.Range("G6").Value = "方法 「{0}」 可能缺少了必需的參數值：{1}。"
$ws.Range("H6").Value = "方法 `"{0}`" 可能缺少了必需的参数值：{1}。"

# Row 6 uses the sheet's (custom) default row height, same as every other
# data row.
$ws.Rows("6").RowHeight = 28.3

# Match the author's final selection/cursor position.
$ws.Range("G6").Select()
